# Updated symbol list on Fri Dec 16 05:38:49 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores numeric-looking values as plain TEXT
# (t="inlineStr" in the original workbook), so we must make sure Excel
# does not silently convert our new values into real numbers when we
# write them back. We do that by temporarily formatting each target
# cell as Text ("@") before assigning the value, then restoring the
# cell style back to Normal/General so we don't leave stray formatting
# behind (the source workbook only changes cell *contents*, not styles).
# NOTE: NumberFormat/Style are applied cell-by-cell (not via a single
# multi-area union range) so every cell is reliably affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceCells = @(
    "D2","D4","D5","D6","D7","D8","D9","D10","D11","D12",
    "D13","D14","D15","D16","D17","D18","D19","D21","D22","D23",
    "D24","D27","D40","D41","D42","D43","D44","D46","D48","D49"
)
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value  = "261.10"
$ws.Range("D4").Value  = "6.209"
$ws.Range("D5").Value  = "0.06153"
$ws.Range("D6").Value  = "6.723"
$ws.Range("D7").Value  = "3.465"
$ws.Range("D8").Value  = "1.342"
$ws.Range("D9").Value  = "0.7984"
$ws.Range("D10").Value = "0.1588"
$ws.Range("D11").Value = "0.08136"
$ws.Range("D12").Value = "0.03496"
$ws.Range("D13").Value = "0.03085"
$ws.Range("D14").Value = "0.09306"
$ws.Range("D15").Value = "3.846"
$ws.Range("D16").Value = "0.001707"
$ws.Range("D17").Value = "0.04790"
$ws.Range("D18").Value = "0.0006142"
$ws.Range("D19").Value = "0.006213"
$ws.Range("D21").Value = "0.004068"
$ws.Range("D22").Value = "0.0001501"
$ws.Range("D23").Value = "3.691"
$ws.Range("D24").Value = "2.214"
$ws.Range("D27").Value = "0.0003203"
$ws.Range("D40").Value = "0.04617"

# Rows 41-43: the three coins rotated one slot (KickToken, BKEXToken,
# CEJI) and got fresh price data.
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.007124"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1119"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.003302"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "0.01024"
$ws.Range("D46").Value = "0.00005943"
$ws.Range("D48").Value = "0.7002"
$ws.Range("D49").Value = "0.1552"

# Restore the cells to the workbook's normal (unformatted) style now
# that the text values are safely in place.
foreach ($ref in $priceCells) {
    $ws.Range($ref).Style = "Normal"
}
